# Rename sheet from "Property1" to "DataNode" to unify the conception of
# DataNode, DataTable, Entity.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Restore the active-cell selection that was captured the last time the
# workbook was saved (bottom-right frozen pane).
$ws.Range("W37").Select() | Out-Null
